$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "73.362.82"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "3.979.46"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "610.52"
$ws.Range("E5").Value = "  +6.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.75"
$ws.Range("E6").Value = "  +10.65%  "
$ws.Range("E7").Value = "  -2.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.789"
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  +8.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.31"
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.25"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "4.620.64"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("D15").Value = "3.978.89"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.32"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.85"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "73.323.67"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.131"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "453.52"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  +5.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.10"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("E24").Value = "  -4.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.23"
$ws.Range("E25").Value = "  -4.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.18"
$ws.Range("E26").Value = "  -2.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.07"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  +0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.54"
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.01"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.97"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").Value = "  +16.11%  "
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "47.94"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.77"
$ws.Range("E36").Value = "  +3.47%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "648.98"
$ws.Range("E37").Value = "  -5.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.431"
$ws.Range("E38").Value = "  -3.77%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.147"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.38"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.26"
$ws.Range("E43").Value = "  +41.27%  "
$ws.Range("E44").Value = "  -2.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("E45").Value = "  -5.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.14"
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("E47").Value = "  -2.76%  "
$ws.Range("E48").Value = "  +8.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.46"
$ws.Range("E49").Value = "  +3.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("E50").Value = "  -4.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.00"
$ws.Range("E51").Value = "  -3.72%  "
